# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape of the source data.

$wb = $excel.ActiveWorkbook

# Row -> new value map for the "展览" sheet (column F), keyed by row number.
$exhibitionUpdates = @{
    3  = 331
    4  = 1443
    5  = 8637
    6  = 81
    7  = 485
    8  = 633
    11 = 3487
    15 = 1074
    19 = 182
    20 = 2216
    21 = 35
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $ws1.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value map for the "全部类型" sheet (column F), keyed by row number.
# (Same underlying records, shifted down by one row because the combined
# sheet also includes the single row from "演出".)
$allTypesUpdates = @{
    3  = 331
    4  = 1443
    5  = 8637
    6  = 81
    7  = 485
    8  = 633
    11 = 3487
    15 = 1074
    19 = 182
    20 = 2216
    22 = 35
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $ws4.Range("F$row").Value = $allTypesUpdates[$row]
}
